$wb = $excel.ActiveWorkbook

# The Application sheet keeps its data, but the cursor moved to B2 and it is
# no longer the active tab once the new Component sheet is added.
$app = $wb.Worksheets.Item("Application")
$app.Range("B2").Select()

# Add the new "Component" sheet after the last existing sheet ("Missing app")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$comp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$comp.Name = "Component"

# Header row
$comp.Range("A1").Value = "application.name"
$comp.Range("B1").Value = "component.id"
$comp.Range("C1").Value = "component.name"

# Data row
$comp.Range("A2").Value = "APPLICATION-0001"
$comp.Range("B2").Value = "SUB.0000001"
$comp.Range("C2").Value = "COMPONENT-0001-0001"

# Column widths (approximate match to the imported layout)
$comp.Columns.Item(1).ColumnWidth = 26.02
$comp.Columns.Item(2).ColumnWidth = 17.23
$comp.Columns.Item(3).ColumnWidth = 21.16

# Leave the cursor one row below the data, ready for the next import row,
# and make this newly-imported sheet the active tab/view.
$comp.Range("C3").Select()
